# Werkblad toevoegen voor 08917
#
# The new sheet "08917" is a duplicate of "08900" (same two-column
# Winkel/Aantal layout, same AutoFilter range/sort state), inserted right
# after "08900" and registered with its own hidden _FilterDatabase
# defined name - mirroring how the existing "08899"/"08900" sheets are set
# up.

$wb = $excel.ActiveWorkbook

# Duplicate the "08900" worksheet (keeps its columns, data, styles and
# AutoFilter/sortState) and drop the copy right after it.
$source = $wb.Worksheets.Item("08900")
$source.Copy($null, $source)

# The copy becomes the new last sheet; rename it to "08917".
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "08917"

# Register the sheet-scoped AutoFilter defined name (_xlnm._FilterDatabase)
# for the new sheet, same as the other two worksheets have.
$newSheet.Names.Add("_xlnm._FilterDatabase", "='08917'!`$A`$1:`$B`$7")
$filterName = $wb.Names.Item($wb.Names.Count)
$filterName.Visible = $false

# Restore the original active sheet/selection (the new sheet otherwise
# ends up as the active/selected tab just because it was created last).
$wb.Worksheets.Item("08899").Activate()
